# Weekly data update: insert two new price records (rows 84-85) for
# "Vega Modelo de Temuco" / Arandano (blue), pushing the existing
# rows 84-103 down to 86-105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 84, shifting rows 84:103 -> 86:105.
$ws.Rows("84:85").Insert()

# --- New row 84 ---
$ws.Cells.Item(84, 1).Value  = 10
$ws.Cells.Item(84, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(84, 3).Value  = "La Araucanía"
$ws.Cells.Item(84, 4).Value  = 44889
$ws.Cells.Item(84, 5).Value  = 9
$ws.Cells.Item(84, 6).Value  = "Fruta"
$ws.Cells.Item(84, 7).Value  = 100101
$ws.Cells.Item(84, 8).Value  = "Berries"
$ws.Cells.Item(84, 9).Value  = 100101001
$ws.Cells.Item(84, 10).Value = "Arándano (blue)"
$ws.Cells.Item(84, 11).Value = "Sin especificar"
$ws.Cells.Item(84, 12).Value = "Primera"
$ws.Cells.Item(84, 13).Value = 200
$ws.Cells.Item(84, 14).Value = 3000
$ws.Cells.Item(84, 15).Value = 3000
$ws.Cells.Item(84, 16).Value = 3000
$ws.Cells.Item(84, 17).Value = "$/kilo"
$ws.Cells.Item(84, 18).Value = "Región Metropolitana"
$ws.Cells.Item(84, 19).Value = 3000
$ws.Cells.Item(84, 20).Value = 1

# --- New row 85 ---
$ws.Cells.Item(85, 1).Value  = 10
$ws.Cells.Item(85, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(85, 3).Value  = "La Araucanía"
$ws.Cells.Item(85, 4).Value  = 44889
$ws.Cells.Item(85, 5).Value  = 9
$ws.Cells.Item(85, 6).Value  = "Fruta"
$ws.Cells.Item(85, 7).Value  = 100101
$ws.Cells.Item(85, 8).Value  = "Berries"
$ws.Cells.Item(85, 9).Value  = 100101001
$ws.Cells.Item(85, 10).Value = "Arándano (blue)"
$ws.Cells.Item(85, 11).Value = "Sin especificar"
$ws.Cells.Item(85, 12).Value = "Primera"
$ws.Cells.Item(85, 13).Value = 250
$ws.Cells.Item(85, 14).Value = 2500
$ws.Cells.Item(85, 15).Value = 2500
$ws.Cells.Item(85, 16).Value = 2500
$ws.Cells.Item(85, 17).Value = "$/kilo"
$ws.Cells.Item(85, 18).Value = "Región del Maule"
$ws.Cells.Item(85, 19).Value = 2500
$ws.Cells.Item(85, 20).Value = 1
